# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price/profit updates to the Seraph_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1518.25
$ws.Range("I55").Value = 526.7273
$ws.Range("J55").Value = 3699.6
$ws.Range("K55").Value = 526.7273
$ws.Range("L55").Value = 3699.6
$ws.Range("M55").Value = -312.7273
$ws.Range("N55").Value = -4127.6
$ws.Range("H58").Value = 3994
$ws.Range("I58").Value = 307.66666
$ws.Range("J58").Value = 5222.778
$ws.Range("K58").Value = 922.9999799999999
$ws.Range("L58").Value = 15668.334
$ws.Range("M58").Value = -772.9999799999999
$ws.Range("N58").Value = -15968.334
$ws.Range("H97").Value = 4460
$ws.Range("J97").Value = 500
$ws.Range("L97").Value = 1500
$ws.Range("N97").Value = -2492
$ws.Range("H106").Value = 32165.889
$ws.Range("I106").Value = 33356.855
$ws.Range("K106").Value = 33356.855
$ws.Range("M106").Value = -32725.855
$ws.Range("H115").Value = 4035.4
$ws.Range("I115").Value = 4035.4
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 12106.2
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -10539.2
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 5999.75
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 1113
$ws.Range("I132").Value = 806.62067
$ws.Range("K132").Value = 2419.86201
$ws.Range("M132").Value = 110.1379900000002
$ws.Range("H135").Value = 956.2727
$ws.Range("I135").Value = 956.2727
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8606.454299999999
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6071.454299999999
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2237
$ws.Range("I2").Value = 2237
$ws.Range("K2").Value = 2237
$ws.Range("M2").Value = -2124
$ws.Range("H74").Value = 1826.9166
$ws.Range("I74").Value = 869.2222
$ws.Range("K74").Value = 869.2222
$ws.Range("M74").Value = 4.777799999999957
$ws.Range("H77").Value = 1826.9166
$ws.Range("I77").Value = 869.2222
$ws.Range("K77").Value = 4346.111
$ws.Range("M77").Value = 21.88900000000012
$ws.Range("H116").Value = 2237
$ws.Range("I116").Value = 2237
$ws.Range("K116").Value = 2237
$ws.Range("M116").Value = 57
$ws.Range("H122").Value = 591179.25
$ws.Range("I122").Value = 669303.1
$ws.Range("K122").Value = 2007909.3
$ws.Range("M122").Value = -2005459.3
$ws.Range("H132").Value = 1096.9166
$ws.Range("I132").Value = 1096.9166
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3290.7498
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -760.7498000000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2237
$ws.Range("I3").Value = 2237
$ws.Range("K3").Value = 2237
$ws.Range("M3").Value = -2123
$ws.Range("H134").Value = 2270
$ws.Range("I134").Value = 2270
$ws.Range("K134").Value = 6810
$ws.Range("M134").Value = -4275

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 6731.6665
$ws.Range("I3").Value = 202
$ws.Range("J3").Value = 9996.5
$ws.Range("K3").Value = 202
$ws.Range("L3").Value = 9996.5
$ws.Range("M3").Value = -89
$ws.Range("N3").Value = -10222.5
$ws.Range("H7").Value = 108.71429
$ws.Range("I7").Value = 58.764706
$ws.Range("K7").Value = 58.764706
$ws.Range("M7").Value = 54.235294
$ws.Range("H31").Value = 3332.5386
$ws.Range("I31").Value = 2865.2727
$ws.Range("K31").Value = 2865.2727
$ws.Range("M31").Value = -2570.2727
$ws.Range("H34").Value = 3332.5386
$ws.Range("I34").Value = 2865.2727
$ws.Range("K34").Value = 2865.2727
$ws.Range("M34").Value = -2663.2727
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H58").Value = 2223.3462
$ws.Range("I58").Value = 1189.55
$ws.Range("K58").Value = 1189.55
$ws.Range("M58").Value = -986.55
$ws.Range("H60").Value = 14767.5
$ws.Range("I60").Value = 10366.3125
$ws.Range("J60").Value = 49977
$ws.Range("K60").Value = 10366.3125
$ws.Range("L60").Value = 49977
$ws.Range("M60").Value = -9855.3125
$ws.Range("N60").Value = -50999
$ws.Range("H86").Value = 6565.4287
$ws.Range("I86").Value = 4190
$ws.Range("K86").Value = 4190
$ws.Range("M86").Value = -3067
$ws.Range("H89").Value = 6565.4287
$ws.Range("I89").Value = 4190
$ws.Range("K89").Value = 20950
$ws.Range("M89").Value = -15334
$ws.Range("H105").Value = 1720.6666
$ws.Range("I105").Value = 1500
$ws.Range("K105").Value = 1500
$ws.Range("M105").Value = 247
$ws.Range("H107").Value = 1463.1666
$ws.Range("I107").Value = 719.75
$ws.Range("K107").Value = 719.75
$ws.Range("M107").Value = 1200.25
$ws.Range("H122").Value = 3787.0908
$ws.Range("I122").Value = 4382.375
$ws.Range("J122").Value = 2199.6667
$ws.Range("K122").Value = 13147.125
$ws.Range("L122").Value = 6599.000100000001
$ws.Range("M122").Value = -10697.125
$ws.Range("N122").Value = -11499.0001
$ws.Range("H134").Value = 2445.1785
$ws.Range("I134").Value = 2273.3809
$ws.Range("K134").Value = 6820.1427
$ws.Range("M134").Value = -4285.1427
$ws.Range("H136").Value = 2223.3462
$ws.Range("I136").Value = 1189.55
$ws.Range("K136").Value = 3568.65
$ws.Range("M136").Value = -1018.65

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1000
$ws.Range("I59").Value = 1000
$ws.Range("J59").Value = 1000
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 3000
$ws.Range("M59").Value = -2460
$ws.Range("N59").Value = -4080
$ws.Range("H68").Value = 1853.8462
$ws.Range("J68").Value = 3550
$ws.Range("L68").Value = 10650
$ws.Range("N68").Value = -12272
$ws.Range("H71").Value = 1853.8462
$ws.Range("J71").Value = 3550
$ws.Range("L71").Value = 31950
$ws.Range("N71").Value = -40062

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 4350
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 4350
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 4350
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -5408
$ws.Range("H41").Value = 780.2
$ws.Range("I41").Value = 633.6667
$ws.Range("J41").Value = 1000
$ws.Range("K41").Value = 633.6667
$ws.Range("L41").Value = 1000
$ws.Range("M41").Value = -278.6667
$ws.Range("N41").Value = -1710
$ws.Range("H80").Value = 9714.5625
$ws.Range("I80").Value = 3809.3333
$ws.Range("K80").Value = 3809.3333
$ws.Range("M80").Value = -2811.3333
$ws.Range("H83").Value = 9714.5625
$ws.Range("I83").Value = 3809.3333
$ws.Range("K83").Value = 19046.6665
$ws.Range("M83").Value = -14054.6665
$ws.Range("H97").Value = 479.5
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 12224.375
$ws.Range("I25").Value = 11680
$ws.Range("J25").Value = 13131.667
$ws.Range("K25").Value = 11680
$ws.Range("L25").Value = 13131.667
$ws.Range("M25").Value = -11450
$ws.Range("N25").Value = -13591.667
$ws.Range("H40").Value = 3166.3333
$ws.Range("I40").Value = 3045.182
$ws.Range("J40").Value = 4499
$ws.Range("K40").Value = 3045.182
$ws.Range("L40").Value = 4499
$ws.Range("M40").Value = -2909.182
$ws.Range("N40").Value = -4771
$ws.Range("H61").Value = 3954.8
$ws.Range("I61").Value = 3942.25
$ws.Range("J61").Value = 4005
$ws.Range("K61").Value = 3942.25
$ws.Range("L61").Value = 4005
$ws.Range("M61").Value = -3740.25
$ws.Range("N61").Value = -4409
$ws.Range("H88").Value = 25000
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 25000
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H113").Value = 3954.8
$ws.Range("I113").Value = 3942.25
$ws.Range("J113").Value = 4005
$ws.Range("K113").Value = 3942.25
$ws.Range("L113").Value = 4005
$ws.Range("M113").Value = -1772.25
$ws.Range("N113").Value = -8345
$ws.Range("H132").Value = 4109.706
$ws.Range("I132").Value = 3991.5625
$ws.Range("K132").Value = 11974.6875
$ws.Range("M132").Value = -9444.6875
$ws.Range("H136").Value = 3349.1875
$ws.Range("I136").Value = 2979.9333
$ws.Range("K136").Value = 8939.7999
$ws.Range("M136").Value = -6389.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 790.7222
$ws.Range("I113").Value = 701
$ws.Range("J113").Value = 880.44446
$ws.Range("K113").Value = 2103
$ws.Range("L113").Value = 2641.33338
$ws.Range("M113").Value = 67
$ws.Range("N113").Value = -6981.33338
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 7968.2856
$ws.Range("I132").Value = 1596
$ws.Range("K132").Value = 4788
$ws.Range("M132").Value = -2258
$ws.Range("H136").Value = 1750.9166
$ws.Range("I136").Value = 1750.9166
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5252.7498
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2702.7498
$ws.Range("N136").ClearContents()
